$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): B11 4 -> 5, C11 -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): B12 72 -> 90, C12 -3 -> -3.6, E12 "69/112" -> "86.4/140"
$ws.Range("B12").Value = 90
$ws.Range("C12").Value = -3.6
$ws.Range("E12").Value = "86.4/140"
